# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Update case statistics for the affected countries
# - Arabia Saudita overtakes Pakistan (rows 16/17 swap)
# - Ghana overtakes Kirguistan (rows 54/55 swap)
# - Principado de Andorra overtakes Republica del Chad (rows 150/151 swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 15:01"

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 4814440
$ws.Cells.Item(4, 3).Value = 793
$ws.Cells.Item(4, 4).Value = 2380561
$ws.Cells.Item(4, 5).Value = 2275504

# India (row 6)
$ws.Cells.Item(6, 2).Value = 1812770
$ws.Cells.Item(6, 3).Value = 8068
$ws.Cells.Item(6, 4).Value = 1190736
$ws.Cells.Item(6, 5).Value = 583785
$ws.Cells.Item(6, 7).Value = 88
$ws.Cells.Item(6, 8).Value = 38249

# Row 16: now Arabia Saudita (passes Pakistan)
$ws.Cells.Item(16, 1).Value = "Arabia Saudita"
$ws.Cells.Item(16, 2).Value = 280093
$ws.Cells.Item(16, 3).Value = 1258
$ws.Cells.Item(16, 4).Value = 242055
$ws.Cells.Item(16, 5).Value = 35089
$ws.Cells.Item(16, 7).Value = 32
$ws.Cells.Item(16, 8).Value = 2949

# Row 17: now Pakistan (pushed down one spot, stats unchanged)
$ws.Cells.Item(17, 1).Value = "Pakistan"
$ws.Cells.Item(17, 2).Value = 280029
$ws.Cells.Item(17, 3).Value = 331
$ws.Cells.Item(17, 4).Value = 248873
$ws.Cells.Item(17, 5).Value = 25172
$ws.Cells.Item(17, 7).Value = 8
$ws.Cells.Item(17, 8).Value = 5984

# Irak (row 24)
$ws.Cells.Item(24, 2).Value = 131886
$ws.Cells.Item(24, 3).Value = 2735
$ws.Cells.Item(24, 4).Value = 94111
$ws.Cells.Item(24, 5).Value = 32841
$ws.Cells.Item(24, 7).Value = 66
$ws.Cells.Item(24, 8).Value = 4934

# Paises Bajos (row 44)
$ws.Cells.Item(44, 2).Value = 55470
$ws.Cells.Item(44, 3).Value = 372

# Portugal (row 47)
$ws.Cells.Item(47, 2).Value = 51569
$ws.Cells.Item(47, 3).Value = 106
$ws.Cells.Item(47, 4).Value = 37111
$ws.Cells.Item(47, 5).Value = 12720

# Row 54: now Ghana (passes Kirguistan)
$ws.Cells.Item(54, 1).Value = "Ghana"
$ws.Cells.Item(54, 2).Value = 37812
$ws.Cells.Item(54, 3).Value = 798
$ws.Cells.Item(54, 4).Value = 34313
$ws.Cells.Item(54, 5).Value = 3308
$ws.Cells.Item(54, 7).Value = 9
$ws.Cells.Item(54, 8).Value = 191

# Row 55: now Kirguistan (pushed down one spot, stats unchanged)
$ws.Cells.Item(55, 1).Value = "Kirguistan"
$ws.Cells.Item(55, 2).Value = 37129
$ws.Cells.Item(55, 3).Value = 410
$ws.Cells.Item(55, 4).Value = 27927
$ws.Cells.Item(55, 5).Value = 7782
$ws.Cells.Item(55, 7).Value = 11
$ws.Cells.Item(55, 8).Value = 1420

# Kenia (row 66)
$ws.Cells.Item(66, 2).Value = 22597
$ws.Cells.Item(66, 3).Value = 544
$ws.Cells.Item(66, 4).Value = 8740
$ws.Cells.Item(66, 5).Value = 13475
$ws.Cells.Item(66, 7).Value = 13
$ws.Cells.Item(66, 8).Value = 382

# Madagascar (row 83)
$ws.Cells.Item(83, 2).Value = 11660
$ws.Cells.Item(83, 3).Value = 132
$ws.Cells.Item(83, 4).Value = 8825
$ws.Cells.Item(83, 5).Value = 2717
$ws.Cells.Item(83, 7).Value = 4
$ws.Cells.Item(83, 8).Value = 118

# Republica de Macedonia (row 84)
$ws.Cells.Item(84, 2).Value = 11128
$ws.Cells.Item(84, 3).Value = 74
$ws.Cells.Item(84, 4).Value = 6972
$ws.Cells.Item(84, 5).Value = 3656
$ws.Cells.Item(84, 7).Value = 3
$ws.Cells.Item(84, 8).Value = 500

# Noruega (row 86)
$ws.Cells.Item(86, 5).Value = 260
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 256

# Croacia (row 100)
$ws.Cells.Item(100, 2).Value = 5294
$ws.Cells.Item(100, 3).Value = 34
$ws.Cells.Item(100, 4).Value = 4438
$ws.Cells.Item(100, 5).Value = 703
$ws.Cells.Item(100, 7).Value = 4
$ws.Cells.Item(100, 8).Value = 153

# Estonia (row 132)
$ws.Cells.Item(132, 2).Value = 1915
$ws.Cells.Item(132, 3).Value = 8
$ws.Cells.Item(132, 5).Value = 80

# Letonia (row 144)
$ws.Cells.Item(144, 2).Value = 1195
$ws.Cells.Item(144, 3).Value = 13
$ws.Cells.Item(144, 4).Value = 1070
$ws.Cells.Item(144, 5).Value = 120
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 5

# Row 150: now Principado de Andorra (passes Republica del Chad)
$ws.Cells.Item(150, 1).Value = "Principado de Andorra"
$ws.Cells.Item(150, 2).Value = 937
$ws.Cells.Item(150, 3).Value = 12
$ws.Cells.Item(150, 4).Value = 821
$ws.Cells.Item(150, 5).Value = 64
$ws.Cells.Item(150, 8).Value = 52

# Row 151: now Republica del Chad (pushed down one spot, stats unchanged)
$ws.Cells.Item(151, 1).Value = "Republica del Chad"
$ws.Cells.Item(151, 2).Value = 936
$ws.Cells.Item(151, 4).Value = 813
$ws.Cells.Item(151, 5).Value = 48
$ws.Cells.Item(151, 8).Value = 75
